$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume symbol data (price refresh snapshot)
$updates = @{
    "D2" = "248.46"
    "D3" = "24.03"
    "D4" = "5.864"
    "D5" = "0.05882"
    "D6" = "3.431"
    "D7" = "6.517"
    "D8" = "1.330"
    "D9" = "0.7974"
    "D10" = "0.1484"
    "D11" = "0.07725"
    "D12" = "0.03303"
    "D13" = "0.03011"
    "D14" = "0.09232"
    "D15" = "3.571"
    "D16" = "0.001663"
    "D17" = "0.04765"
    "D18" = "0.0006034"
    "E18" = "17OneONEWorstin24h"
    "D19" = "0.006225"
    "D20" = "0.005536"
    "D21" = "0.001070"
    "D22" = "0.0001502"
    "D23" = "3.727"
    "D26" = "0.1253"
    "D27" = "0.0006247"
    "D40" = "0.04367"
    "D41" = "0.007041"
    "D43" = "0.003410"
    "D44" = "0.009648"
    "D45" = "0.002462"
    "E45" = "44ACDXExchangeACXT"
    "D46" = "0.00005901"
    "D48" = "0.9906"
    "D49" = "0.1076"
    "E49" = "48BOLOBOLOBestin24h"
    "D51" = "0.01011"
}

foreach ($cellRef in $updates.Keys) {
    $newValue = $updates[$cellRef]
    $range = $ws.Range($cellRef)
    $col = $cellRef.Substring(0, 1)
    if ($col -eq "D") {
        # Price column stores numeric-looking values as text; force text format
        # so Excel does not silently coerce the literal to a Number cell,
        # then restore the default cell style so formatting is unaffected.
        $range.NumberFormat = "@"
        $range.Value = $newValue
        $range.Style = "Normal"
    } else {
        $range.Value = $newValue
    }
}
